# B6-PowerPoint.pptx edit:
#  1) Re-point the three "Table_0" tables (slides 14, 15, 16) from the
#     custom table style {2E90BAC3-...} to the built-in table style
#     {04381D02-8A66-493A-AE64-85516CC5F97A}.
#  2) Re-apply the presentation's theme colour scheme so the deck uses
#     the default "Office" palette instead of the "Integral / Red
#     Violet" palette that was previously in force.

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newTableStyle = "{04381D02-8A66-493A-AE64-85516CC5F97A}"
foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    $shape = $slide.Shapes.Item(1)
    if ($shape.HasTable) {
        $shape.Table.ApplyStyle($newTableStyle)
    }
}

# --- 2. Theme colour scheme -------------------------------------------
# Office default theme colours (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
# expressed as COM BGR-packed RGB integers.
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$slide1 = $p.Slides.Item(1)
$themeColors = $slide1.ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
